$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 with the new shared string, copying B1's style (bordered, bold, centered)
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Other parameter values:"

# Set the new column's width
$ws.Columns.Item(8).ColumnWidth = 27

# Update A10's style: same as B1 (bold/center/top) but without the border
$ws.Range("B1").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Borders.LineStyle = -4142

# Reset selection back to A1 (default / top-left)
$ws.Range("A1").Select()
